$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")
$ws.Activate()

$ws.Columns("N").Insert()
$ws.Columns("N").ColumnWidth = $ws.Columns("M").ColumnWidth

$ws.Range("M13").Select() | Out-Null
